# CP-159: Fix Student template data (SSID/AlternateSSID/GradeLevelWhenAssessed)
# for rows 4-7 so the sample rows reflect students in grades 07 and 11 rather
# than duplicating grade 03, and update the active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 -> grade 07, first student
$ws.Range("H4").Value = "IRP79990001"
$ws.Range("I4").Value = "AIRP79990001"
$ws.Range("J4").Value = "07"

# Row 5 -> grade 07, second student
$ws.Range("H5").Value = "IRP79990002"
$ws.Range("I5").Value = "AIRP79990002"
$ws.Range("J5").Value = "07"

# Row 6 -> grade 11, first student
$ws.Range("H6").Value = "IRP119990001"
$ws.Range("I6").Value = "AIRP119990001"
$ws.Range("J6").Value = "11"

# Row 7 -> grade 11, second student
$ws.Range("H7").Value = "IRP119990002"
$ws.Range("I7").Value = "AIRP119990002"
$ws.Range("J7").Value = "11"

# Update the sheet's active cell / selection to I6 (matches saved view state)
$ws.Range("I6").Select()
